# Add the missing "OK" button row (code A0 / 5F) to the Acer TV (1363)
# button-code table on the first worksheet, just above the "SUB-T" entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# Insert a new blank row at row 30 - this pushes the existing rows 30:44
# down to 31:45 (formats/styles are carried down from the row above, which
# already matches what we want for the new data row).
[void]$ws.Rows.Item(30).Insert()

# Fill in the new row with the OK button's name/data.
$ws.Range("A30").Value = "OK"
$ws.Range("B30").Value = "A0"
$ws.Range("C30").Value = "5F"

# Match the saved view/selection state from the edit.
$excel.ActiveWindow.ScrollRow = 21
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("A31").Select()
